# Applies the "Part design and functions" / "Tasks" section edits:
#   - appends descriptive sentences to each component bullet
#   - appends a colon (and in some cases descriptive text) to each task bullet
#   - drops the stray <w:lastRenderedPageBreak/> before "Tasks"

$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd("`r") -eq $text) {
            return $p
        }
    }
    return $null
}

# --- Part design and functions ---

$p = Get-ParaByText $d "Instructional SRAM"
$p.Range.InsertAfter(":  Designed with if statement. When certain conditions are met it turns into a case when statement. ")
$p.Range.InsertAfter("It is read only memory so you can't change it internally. ")

$p = Get-ParaByText $d "Control Unit"
$p.Range.InsertAfter(": Reads instruction SRAM in machine code and decodes it. Turns registers off and on. The ALU is always on. It will write the output of the ALU to the Data SRAM. ")

$p = Get-ParaByText $d "Program Counter"
$p.Range.InsertAfter(": ")
$p.Range.InsertAfter("Its")
$p.Range.InsertAfter(" a register that increments")

$p = Get-ParaByText $d "ALU"
$p.Range.InsertAfter(":")
$p.Range.InsertAfter(" 32 2-bit full adder with 2-bit XOR gates on the inputs. It does addition and subtraction only. ")

$p = Get-ParaByText $d "Registers"
$p.Range.InsertAfter(":")
$p.Range.InsertAfter(" 32 D flip flops with 8-1 multiplexor. ")

$p = Get-ParaByText $d "Data SRAM"
$p.Range.InsertAfter(":")
$p.Range.InsertAfter(" ")

# --- Drop the stray lastRenderedPageBreak before "Tasks" ---

$p = Get-ParaByText $d "Tasks"
$p.Range.Text = "Tasks"

# --- Tasks ---

$p = Get-ParaByText $d "Task 1"
$p.Range.InsertAfter(": it works. it does the do stuff yes")

$p = Get-ParaByText $d "Task 2"
$p.Range.InsertAfter(":")

$p = Get-ParaByText $d "Task 3"
$p.Range.InsertAfter(":")

$p = Get-ParaByText $d "Task 4"
$p.Range.InsertAfter(":")
